$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column E width (approximates the 60.83203125 target width; COM rounds to 1/6-character granularity)
$ws.Columns.Item(5).ColumnWidth = 60

# Header row
$ws.Range("A1").Value = 'Vendor'
$ws.Range("B1").Value = 'Website'
$ws.Range("C1").Value = 'Welcome Message Evaluation'
$ws.Range("D1").Value = 'Was Detected'
$ws.Range("E1").Value = 'Screenshot Path'

# Row 2
$ws.Range("A2").Value = 'Tidio'
$ws.Range("B2").Value = 'https://www.tidio.com/'
$ws.Range("C2").Value = 'The chatbot widget is visible. The welcome message effectively introduces the functionalities by listing clear and relevant options for user assistance. The tone is friendly and inviting, offering help with specific inquiries while also providing a chat option with AI support when the team is offline. 
Score: 90/100'
$ws.Range("D2").Value = 'Yes'
$ws.Range("E2").Value = 'C:\Users\elmun\OneDrive\Escritorio\Zoros\Trabajo\Proyecto OO2\Chatbot-Vendor-Detector\vendor_screenshots\tidio_1743546415881.jpg'

# Row 3
$ws.Range("A3").Value = 'LiveChat'
$ws.Range("B3").Value = 'https://www.LiveChat.com/'
$ws.Range("C3").Value = 'The chatbot widget is visible. The welcome message effectively introduces the main functionality by highlighting engagement tools and invites interaction. 
Score: 85'
$ws.Range("D3").Value = 'Yes'
$ws.Range("E3").Value = 'C:\Users\elmun\OneDrive\Escritorio\Zoros\Trabajo\Proyecto OO2\Chatbot-Vendor-Detector\vendor_screenshots\livechat_1743546440780.jpg'

# Row 4
$ws.Range("A4").Value = 'HubSpot'
$ws.Range("B4").Value = 'https://www.hubspot.com/'
$ws.Range("C4").Value = 'The chatbot widget is visible. 
The welcome message effectively presents clear options for user engagement, such as "Get free training," "Get started free," "Book a demo," and "Chat with sales." These options introduce the chatbot''s functionalities well and guide user interaction effectively.
Score: 85/100'
$ws.Range("D4").Value = 'Yes'
$ws.Range("E4").Value = 'C:\Users\elmun\OneDrive\Escritorio\Zoros\Trabajo\Proyecto OO2\Chatbot-Vendor-Detector\vendor_screenshots\hubspot_1743546472266.jpg'

# Row 5
$ws.Range("A5").Value = 'ChatBot.com'
$ws.Range("B5").Value = 'https://www.chatbot.com/'
$ws.Range("C5").Value = 'The chatbot widget is visible. The welcome message is friendly and interactive, offering clear navigation options to guide users through different functionalities like building a chatbot, using it, asking questions, or just browsing. The message effectively introduces the chatbot''s capabilities. Score: 90/100.'
$ws.Range("D5").Value = 'Yes'
$ws.Range("E5").Value = 'C:\Users\elmun\OneDrive\Escritorio\Zoros\Trabajo\Proyecto OO2\Chatbot-Vendor-Detector\vendor_screenshots\chatbot.com_1743546495103.jpg'

# Row 6
$ws.Range("A6").Value = 'Trengo'
$ws.Range("B6").Value = 'https://trengo.com/'
$ws.Range("C6").Value = 'No chatbot widget found for analysis.'
$ws.Range("D6").Value = 'No'
$ws.Range("E6").Value = 'C:\Users\elmun\OneDrive\Escritorio\Zoros\Trabajo\Proyecto OO2\Chatbot-Vendor-Detector\vendor_screenshots\trengo_1743546510074.jpg'

# Row 7
$ws.Range("A7").Value = 'Verloop'
$ws.Range("B7").Value = 'https://www.verloop.io/'
$ws.Range("C7").Value = 'No chatbot widget found for analysis.'
$ws.Range("D7").Value = 'No'
$ws.Range("E7").Value = 'C:\Users\elmun\OneDrive\Escritorio\Zoros\Trabajo\Proyecto OO2\Chatbot-Vendor-Detector\vendor_screenshots\verloop_1743546527205.jpg'

# Row 8
$ws.Range("A8").Value = 'kommunicate'
$ws.Range("B8").Value = 'https://www.kommunicate.io/'
$ws.Range("C8").Value = 'The chatbot widget is visible, but it doesn''t provide any welcome messages introducing functionalities or tutorials. Score: 20.'
$ws.Range("D8").Value = 'Yes'
$ws.Range("E8").Value = 'C:\Users\elmun\OneDrive\Escritorio\Zoros\Trabajo\Proyecto OO2\Chatbot-Vendor-Detector\vendor_screenshots\kommunicate_1743546549118.jpg'
